# Nudge the "5-Point Star 12" shape on the Conclusion slide (slide 9) to its
# new position.
#   old OOXML offset: <a:off x="7969063" y="4112895"/>
#   new OOXML offset: <a:off x="7956998" y="4514215"/>
# PowerPoint COM exposes Shape.Left/Shape.Top in points (1 pt = 12700 EMU),
# so the EMU offsets from the target OOXML are converted accordingly.

$p = $ppt.ActivePresentation

$star = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.Name -eq "5-Point Star 12" -and [Math]::Round($shape.Left * 12700) -eq 7969063) {
            $star = $shape
        }
    }
    if ($star -ne $null) { break }
}

$star.Left = 7956998 / 12700
$star.Top = 4514215 / 12700
